$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) values round-trip as text, not auto-converted numbers,
# matching the source data (inline strings in the original workbook).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.347.20'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '1.867.27'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '234.51'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '0.4706'
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("D8").Value = '0.2871'
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("D9").Value = '0.06566'
$ws.Range("E9").Value = '  +0.49%  '
$ws.Range("D10").Value = '21.47'
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").Value = '0.07874'
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("D12").Value = '96.87'
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("D13").Value = '1.867.88'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '0.6921'
$ws.Range("E14").Value = '  +1.72%  '
$ws.Range("D15").Value = '5.105'
$ws.Range("D16").Value = '268.02'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = '30.348.51'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").Value = '13.96'
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").Value = '0.000007667'
$ws.Range("E19").Value = '  +3.68%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").Value = '2.120.33'
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '5.237'
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("D24").Value = '6.188'
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").Value = '9.397'
$ws.Range("E25").Value = '  +2.21%  '
$ws.Range("D26").Value = '167.58'
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").Value = '18.86'
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = '1.946'
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.361'
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.09919'
$ws.Range("E30").Value = '  +0.96%  '
$ws.Range("D31").Value = '4.386'
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("D32").Value = '1.461'
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("D33").Value = '4.056'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = '0.04740'
$ws.Range("E34").Value = '  +0.89%  '
$ws.Range("D35").Value = '1.133'
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("D36").Value = '0.7030'
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").Value = '2.720'
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").Value = '0.01871'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = '2.800'
$ws.Range("E39").Value = '  +7.09%  '
$ws.Range("D40").Value = '6.281'
$ws.Range("E40").Value = '  +0.42%  '
$ws.Range("D41").Value = '73.49'
$ws.Range("E41").Value = '  -1.08%  '
$ws.Range("D42").Value = '1.947'
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("D43").Value = '0.8431'
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("D44").Value = '0.4174'
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").Value = '103.21'
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").Value = '981.61'
$ws.Range("E47").Value = '  +2.59%  '
$ws.Range("D48").Value = '7.104'
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("D49").Value = '9.109'
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("D50").Value = '34.49'
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("D51").Value = '0.05676'
$ws.Range("E51").Value = '  +0.31%  '

# Reset the style index back to the default ("Normal") now that the text
# values are committed, so no extra cell-level style attribute lingers.
$ws.Range("D2:D51").Style = "Normal"
